{"js": "// Map of old \"dividend\u00f7divisor=\" text -> new \"dividend\u00f7divisor=\" text,\n// taken from the unified diff of the worksheet table cells.\nconst replacements = {\n  \"64\u00f75=\": \"67\u00f79=\",\n  \"57\u00f72=\": \"45\u00f72=\",\n  \"65\u00f79=\": \"88\u00f79=\",\n  \"62\u00f79=\": \"50\u00f78=\",\n  \"92\u00f74=\": \"23\u00f76=\",\n  \"60\u00f74=\": \"51\u00f77=\",\n  \"78\u00f77=\": \"83\u00f76=\",\n  \"79\u00f77=\": \"38\u00f72=\",\n  \"49\u00f73=\": \"19\u00f72=\",\n  \"78\u00f79=\": \"28\u00f73=\",\n  \"93\u00f78=\": \"81\u00f77=\",\n  \"76\u00f75=\": \"49\u00f77=\",\n  \"45\u00f75=\": \"29\u00f75=\",\n  \"36\u00f78=\": \"97\u00f79=\",\n  \"25\u00f78=\": \"76\u00f72=\",\n  \"66\u00f74=\": \"46\u00f76=\",\n  \"42\u00f74=\": \"32\u00f73=\",\n  \"11\u00f75=\": \"45\u00f78=\",\n  \"67\u00f75=\": \"55\u00f78=\",\n  \"48\u00f78=\": \"78\u00f78=\",\n  \"43\u00f76=\": \"40\u00f79=\",\n  \"68\u00f72=\": \"66\u00f75=\",\n  \"69\u00f72=\": \"52\u00f77=\",\n  \"62\u00f72=\": \"61\u00f79=\",\n  \"29\u00f78=\": \"18\u00f75=\",\n};\n\nconst body = context.document.body;\n\nfor (const oldText of Object.keys(replacements)) {\n  const newText = replacements[oldText];\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (const range of found.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Map of old \"dividend\u00f7divisor=\" text -> new \"dividend\u00f7divisor=\" text,\n# taken from the unified diff of the worksheet table cells.\n$replacements = [ordered]@{\n    \"64\u00f75=\" = \"67\u00f79=\"\n    \"57\u00f72=\" = \"45\u00f72=\"\n    \"65\u00f79=\" = \"88\u00f79=\"\n    \"62\u00f79=\" = \"50\u00f78=\"\n    \"92\u00f74=\" = \"23\u00f76=\"\n    \"60\u00f74=\" = \"51\u00f77=\"\n    \"78\u00f77=\" = \"83\u00f76=\"\n    \"79\u00f77=\" = \"38\u00f72=\"\n    \"49\u00f73=\" = \"19\u00f72=\"\n    \"78\u00f79=\" = \"28\u00f73=\"\n    \"93\u00f78=\" = \"81\u00f77=\"\n    \"76\u00f75=\" = \"49\u00f77=\"\n    \"45\u00f75=\" = \"29\u00f75=\"\n    \"36\u00f78=\" = \"97\u00f79=\"\n    \"25\u00f78=\" = \"76\u00f72=\"\n    \"66\u00f74=\" = \"46\u00f76=\"\n    \"42\u00f74=\" = \"32\u00f73=\"\n    \"11\u00f75=\" = \"45\u00f78=\"\n    \"67\u00f75=\" = \"55\u00f78=\"\n    \"48\u00f78=\" = \"78\u00f78=\"\n    \"43\u00f76=\" = \"40\u00f79=\"\n    \"68\u00f72=\" = \"66\u00f75=\"\n    \"69\u00f72=\" = \"52\u00f77=\"\n    \"62\u00f72=\" = \"61\u00f79=\"\n    \"29\u00f78=\" = \"18\u00f75=\"\n}\n\n$d = $word.ActiveDocument\n\nforeach ($oldText in $replacements.Keys) {\n    $newText = $replacements[$oldText]\n\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Replacement.ClearFormatting()\n\n    # FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,\n    # MatchAllWordForms, Forward, Wrap (1=wdFindContinue), Format,\n    # ReplaceWith, Replace (2=wdReplaceAll)\n    $range.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n"}
